$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.165.34"
$ws.Range("E2").Value = "  +3.59%  "
$ws.Range("D3").Value = "3.063.22"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'560.89"
$ws.Range("E5").Value = "  +3.58%  "
$ws.Range("D6").Value = "'143.72"
$ws.Range("E6").Value = "  +4.47%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "3.064.26"
$ws.Range("E8").Value = "  +2.62%  "
$ws.Range("E9").Value = "  +5.45%  "
$ws.Range("E10").Value = "  +6.64%  "
$ws.Range("D11").Value = "'6.10"
$ws.Range("E11").Value = "  -8.72%  "
$ws.Range("D12").Value = "'0.488"
$ws.Range("E12").Value = "  +10.49%  "
$ws.Range("D13").Value = "'0.0000232"
$ws.Range("E13").Value = "  +6.14%  "
$ws.Range("D14").Value = "'35.58"
$ws.Range("E14").Value = "  +5.66%  "
$ws.Range("D15").Value = "3.563.18"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "64.195.15"
$ws.Range("E16").Value = "  +3.47%  "
$ws.Range("D17").Value = "3.070.26"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("E18").Value = "  +2.81%  "
$ws.Range("D19").Value = "'6.77"
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("D20").Value = "'479.07"
$ws.Range("E20").Value = "  +3.48%  "
$ws.Range("D21").Value = "'13.99"
$ws.Range("E21").Value = "  +5.58%  "
$ws.Range("D22").Value = "'0.682"
$ws.Range("E22").Value = "  +5.43%  "
$ws.Range("E23").Value = "  +5.80%  "
$ws.Range("D24").Value = "'14.32"
$ws.Range("E24").Value = "  +14.99%  "
$ws.Range("D25").Value = "'82.31"
$ws.Range("E25").Value = "  +4.44%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +4.10%  "
$ws.Range("D28").Value = "'8.08"
$ws.Range("E28").Value = "  +7.01%  "
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'26.38"
$ws.Range("E31").Value = "  +4.44%  "
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").Value = "'2.45"
$ws.Range("E33").Value = "  +5.65%  "
$ws.Range("D34").Value = "'5.73"
$ws.Range("E34").Value = "  +4.16%  "
$ws.Range("D35").Value = "'6.27"
$ws.Range("E35").Value = "  +8.54%  "
$ws.Range("D36").Value = "'54.96"
$ws.Range("E36").Value = "  +1.19%  "
$ws.Range("E37").Value = "  +5.53%  "
$ws.Range("D38").Value = "'445.20"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "'0.0813"
$ws.Range("E39").Value = "  +1.59%  "
$ws.Range("E40").Value = "  +12.01%  "
$ws.Range("D41").Value = "2.996.44"
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("D42").Value = "'8.27"
$ws.Range("E42").Value = "  +3.50%  "
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "'27.91"
$ws.Range("E44").Value = "  +4.97%  "
$ws.Range("D45").Value = "'0.263"
$ws.Range("E45").Value = "  +7.46%  "
$ws.Range("D46").Value = "'2.18"
$ws.Range("E46").Value = "  +10.33%  "
$ws.Range("E48").Value = "  +4.95%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").Value = "0.0₃0520"
$ws.Range("E49").Value = "  +5.66%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'118.51"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("D51").Value = "'2.08"
$ws.Range("E51").Value = "  +5.07%  "
